# "Generate Report for Handback" — record a handback-transform failure for
# the d8c9a66c... file: update its Status everywhere it is shown (the
# Overview rollup columns + each language sheet's Status column) and fill
# in the Error Detail (column P) on the zh-cn and de-de sheets with the
# localized failure message. Also widen column P so the longer message is
# readable.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$newStatus = "Handback transform failed"

# Overview sheet: zh-cn / de-de status rollup columns for the d8c9a66c row
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

# Per-language sheets: Status column (C) for the same row
$zhcn.Range("C3").Value = $newStatus
$dede.Range("C3").Value = $newStatus

# Per-language sheets: Error Detail column (P) for the same row
$zhcn.Range("P3").Value = "Handback file name: vzc15cv3.hso is different with handoff file name: d8c9a66c-c030-465d-ace6-d21e6c1d337c.112da40fcddd73bbcc6e012f1fb82fe628ce01e9.zh-cn."
$dede.Range("P3").Value = "Handback file name: vzc15cv3.hso is different with handoff file name: d8c9a66c-c030-465d-ace6-d21e6c1d337c.112da40fcddd73bbcc6e012f1fb82fe628ce01e9.de-de."

# Widen the Error Detail column so the message is readable (stored column
# width of 40 chars; ColumnWidth reports ~0.83 less than the stored width
# due to the default-font padding Excel adds).
$zhcn.Columns.Item(16).ColumnWidth = 39.166666666666664
$dede.Columns.Item(16).ColumnWidth = 39.166666666666664
